$d = $word.ActiveDocument

# Locate the text span covering the three runs that together read
# "– organizacija 2 " + "–" + " ponavljanje" == "– organizacija 2 – ponavljanje".
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("– organizacija 2 – ponavljanje", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$spanStart = $findRange.Start
$spanEnd   = $findRange.End

# First new run: "– O" (italic, not bold) -- replaces the leading "– O" of the
# old text (same length, so offsets after it do not move).
$r1 = $d.Range($spanStart, $spanStart + 3)
$r1.Text = "– O"
$r1.Font.Bold = $false
$r1.Font.Italic = $true

# Second new run: "rganizacija 2 – ponavljanje" (same formatting) -- replaces
# the rest of the old text ("rganizacija 2 – ponavljanje", also same length,
# so $spanEnd is still the correct end offset).
$r2 = $d.Range($spanStart + 3, $spanEnd)
$r2.Text = "rganizacija 2 – ponavljanje"
$r2.Font.Bold = $false
$r2.Font.Italic = $true

# Move the "_GoBack" bookmark so that it wraps the newly-formatted text
# (it used to sit right after it, at the very end of the paragraph).
$d.Bookmarks.Add("_GoBack", $d.Range($spanStart, $spanEnd)) | Out-Null

Write-Output $d.Range(0, $spanEnd).Text
